$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 4's formatting down to the new row 5 (keeps the date-style xf
# on A/G reused instead of minting a new numFmt entry), then fill in values.
$ws.Range("A4:I4").Copy()
$ws.Range("A5:I5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(5, 1).Value = 42636.589085648149
$ws.Cells.Item(5, 2).Value = $false
$ws.Cells.Item(5, 3).Value = 9917.16
$ws.Cells.Item(5, 4).Value = 9948
$ws.Cells.Item(5, 5).Value = 19.29
$ws.Cells.Item(5, 6).Value = 19.41
$ws.Cells.Item(5, 7).Value = $true
$ws.Cells.Item(5, 8).Value = 0.62
$ws.Cells.Item(5, 9).Value = $false
